$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count
$newToday = 20260108

for ($r = 2; $r -le $lastRow; $r++) {
    $dVal = $ws.Cells.Item($r, 4).Value2
    $eVal = $ws.Cells.Item($r, 5).Value2
    $fVal = $ws.Cells.Item($r, 6).Value2

    if ($null -eq $dVal -or $null -eq $eVal) {
        continue
    }

    # Row has a malformed 9-digit start date (data-entry typo) - the source
    # process can't compute a remaining-days delta for it, so it is left
    # untouched, same as every other field in that row.
    if ($fVal -gt 99999999) {
        continue
    }

    if ($eVal -eq 1) {
        # Remaining days hit the floor - roll the booking over: reset
        # remaining to the full duration and bump the start date to "today".
        $ws.Cells.Item($r, 5).Value = $dVal
        $ws.Cells.Item($r, 6).Value = $newToday
    } else {
        # One more day has elapsed since the last update.
        $ws.Cells.Item($r, 5).Value = $eVal - 1
    }
}
